$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 663.1
$ws.Range("I2").Value = 396.2
$ws.Range("J2").Value = 930
$ws.Range("K2").Value = 396.2
$ws.Range("L2").Value = 930
$ws.Range("M2").Value = -283.2
$ws.Range("N2").Value = -1156

$ws.Range("H34").Value = 14212.25
$ws.Range("I34").Value = 2283
$ws.Range("K34").Value = 2283
$ws.Range("M34").Value = -2080

$ws.Range("H36").Value = 14212.25
$ws.Range("I36").Value = 2283
$ws.Range("K36").Value = 2283
$ws.Range("M36").Value = -1568

$ws.Range("H70").Value = 1591.4333
$ws.Range("I70").Value = 1061.3572
$ws.Range("J70").Value = 2055.25
$ws.Range("K70").Value = 3184.0716
$ws.Range("L70").Value = 6165.75
$ws.Range("M70").Value = -2914.0716
$ws.Range("N70").Value = -6705.75

$ws.Range("H73").Value = 1591.4333
$ws.Range("I73").Value = 1061.3572
$ws.Range("J73").Value = 2055.25
$ws.Range("K73").Value = 3184.0716
$ws.Range("L73").Value = 6165.75
$ws.Range("M73").Value = -2248.0716
$ws.Range("N73").Value = -8037.75

$ws.Range("H74").Value = 3111.4285
$ws.Range("I74").Value = 2948.3872
$ws.Range("J74").Value = 4375
$ws.Range("K74").Value = 2948.3872
$ws.Range("L74").Value = 4375
$ws.Range("M74").Value = -2012.3872
$ws.Range("N74").Value = -6247

$ws.Range("H77").Value = 3111.4285
$ws.Range("I77").Value = 2948.3872
$ws.Range("J77").Value = 4375
$ws.Range("K77").Value = 14741.936
$ws.Range("L77").Value = 21875
$ws.Range("M77").Value = -10061.936
$ws.Range("N77").Value = -31235

$ws.Range("H137").Value = 1385.4458
$ws.Range("I137").Value = 1445
$ws.Range("J137").Value = 1344.1224
$ws.Range("K137").Value = 4335
$ws.Range("L137").Value = 4032.3672
$ws.Range("M137").Value = -1785
$ws.Range("N137").Value = -9132.367200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4523.643
$ws.Range("I32").Value = 3484.5715
$ws.Range("J32").Value = 11797.143
$ws.Range("K32").Value = 3484.5715
$ws.Range("L32").Value = 11797.143
$ws.Range("M32").Value = -3197.5715
$ws.Range("N32").Value = -12371.143

$ws.Range("H74").Value = 40347.25
$ws.Range("I74").Value = 50899.75
$ws.Range("J74").Value = 13966
$ws.Range("K74").Value = 50899.75
$ws.Range("L74").Value = 13966
$ws.Range("M74").Value = -50025.75
$ws.Range("N74").Value = -15714

$ws.Range("H77").Value = 40347.25
$ws.Range("I77").Value = 50899.75
$ws.Range("J77").Value = 13966
$ws.Range("K77").Value = 254498.75
$ws.Range("L77").Value = 69830
$ws.Range("M77").Value = -250130.75
$ws.Range("N77").Value = -78566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1617.6471
$ws.Range("J20").Value = 1925
$ws.Range("L20").Value = 1925
$ws.Range("N20").Value = -2419

$ws.Range("H44").Value = 19900
$ws.Range("J44").Value = 19900
$ws.Range("L44").Value = 19900
$ws.Range("N44").Value = -20894

$ws.Range("H54").Value = 5777.2144
$ws.Range("I54").Value = 2875.6667
$ws.Range("J54").Value = 11000
$ws.Range("K54").Value = 2875.6667
$ws.Range("L54").Value = 11000
$ws.Range("M54").Value = -2391.6667
$ws.Range("N54").Value = -11968

$ws.Range("H134").Value = 6539.921
$ws.Range("I134").Value = 6270.6
$ws.Range("J134").Value = 7549.875
$ws.Range("K134").Value = 18811.8
$ws.Range("L134").Value = 22649.625
$ws.Range("M134").Value = -16276.8
$ws.Range("N134").Value = -27719.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 42384.6
$ws.Range("I11").Value = 952.5
$ws.Range("J11").Value = 70006
$ws.Range("K11").Value = 952.5
$ws.Range("L11").Value = 70006
$ws.Range("M11").Value = -812.5
$ws.Range("N11").Value = -70286

$ws.Range("H134").Value = 31430282
$ws.Range("I134").Value = 3847831.5
$ws.Range("J134").Value = 111112920
$ws.Range("K134").Value = 11543494.5
$ws.Range("L134").Value = 333338760
$ws.Range("M134").Value = -11540959.5
$ws.Range("N134").Value = -333343830

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3059.889
$ws.Range("J39").Value = 3342.375
$ws.Range("L39").Value = 10027.125
$ws.Range("N39").Value = -10615.125

$ws.Range("H55").Value = 3566.1538
$ws.Range("J55").Value = 3696.6667
$ws.Range("L55").Value = 11090.0001
$ws.Range("N55").Value = -11444.0001

$ws.Range("H68").Value = 1584.0146
$ws.Range("I68").Value = 951.13635
$ws.Range("J68").Value = 1886.6957
$ws.Range("K68").Value = 2853.40905
$ws.Range("L68").Value = 5660.0871
$ws.Range("M68").Value = -2042.40905
$ws.Range("N68").Value = -7282.0871

$ws.Range("H71").Value = 1584.0146
$ws.Range("I71").Value = 951.13635
$ws.Range("J71").Value = 1886.6957
$ws.Range("K71").Value = 8560.227150000001
$ws.Range("L71").Value = 16980.2613
$ws.Range("M71").Value = -4504.227150000001
$ws.Range("N71").Value = -25092.2613

$ws.Range("H107").Value = 218283.52
$ws.Range("I107").Value = 726.7778
$ws.Range("J107").Value = 358141.44
$ws.Range("K107").Value = 2180.3334
$ws.Range("L107").Value = 1074424.32
$ws.Range("M107").Value = -260.3334
$ws.Range("N107").Value = -1078264.32

$ws.Range("H131").Value = 893.45764
$ws.Range("I131").Value = 326.7
$ws.Range("J131").Value = 1009.12244
$ws.Range("K131").Value = 980.0999999999999
$ws.Range("L131").Value = 3027.36732
$ws.Range("M131").Value = 4059.9
$ws.Range("N131").Value = -13107.36732

$ws.Range("H139").Value = 1917.6
$ws.Range("I139").Value = 1496.6666
$ws.Range("K139").Value = 4489.9998
$ws.Range("M139").Value = 650.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8211.857
$ws.Range("J5").Value = 8211.857
$ws.Range("L5").Value = 8211.857
$ws.Range("N5").Value = -8435.857

$ws.Range("H42").Value = 27000
$ws.Range("J42").Value = 27000
$ws.Range("L42").Value = 27000
$ws.Range("N42").Value = -27970

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H115").Value = 27000
$ws.Range("J115").Value = 27000
$ws.Range("L115").Value = 27000
$ws.Range("N115").Value = -29350

$ws.Range("H126").Value = 1765.3077
$ws.Range("I126").Value = 1694.9
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5084.700000000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2614.700000000001
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1187.909
$ws.Range("I16").Value = 1066.6
$ws.Range("J16").Value = 2401
$ws.Range("K16").Value = 1066.6
$ws.Range("L16").Value = 2401
$ws.Range("M16").Value = -896.5999999999999
$ws.Range("N16").Value = -2741

$ws.Range("H46").Value = 792.6070999999999
$ws.Range("I46").Value = 725.0833
$ws.Range("J46").Value = 843.25
$ws.Range("K46").Value = 725.0833
$ws.Range("L46").Value = 843.25
$ws.Range("M46").Value = -537.0833
$ws.Range("N46").Value = -1219.25

$ws.Range("H122").Value = 3422.1428
$ws.Range("I122").Value = 3420
$ws.Range("J122").Value = 3425.625
$ws.Range("K122").Value = 10260
$ws.Range("L122").Value = 10276.875
$ws.Range("M122").Value = -7810
$ws.Range("N122").Value = -15176.875

$ws.Range("H139").Value = 38104.168
$ws.Range("I139").Value = 10650
$ws.Range("J139").Value = 43595
$ws.Range("K139").Value = 10650
$ws.Range("L139").Value = 43595
$ws.Range("M139").Value = -5510
$ws.Range("N139").Value = -53875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 32357.25
$ws.Range("J46").Value = 32357.25
$ws.Range("L46").Value = 32357.25
$ws.Range("N46").Value = -32819.25

$ws.Range("H81").Value = 1683.3334
$ws.Range("J81").Value = 1933
$ws.Range("L81").Value = 3866
$ws.Range("N81").Value = -5988

$ws.Range("H84").Value = 1683.3334
$ws.Range("J84").Value = 1933
$ws.Range("L84").Value = 19330
$ws.Range("N84").Value = -29938

$ws.Range("H126").Value = 32258994
$ws.Range("I126").Value = 536.7
$ws.Range("J126").Value = 90910740
$ws.Range("K126").Value = 1610.1
$ws.Range("L126").Value = 272732220
$ws.Range("M126").Value = 859.8999999999999
$ws.Range("N126").Value = -272737160

$ws.Range("H132").Value = 3299.6584
$ws.Range("I132").Value = 3606.2593
$ws.Range("J132").Value = 2708.3572
$ws.Range("K132").Value = 10818.7779
$ws.Range("L132").Value = 8125.071599999999
$ws.Range("M132").Value = -8288.777900000001
$ws.Range("N132").Value = -13185.0716

$ws.Range("H134").Value = 32357.25
$ws.Range("J134").Value = 32357.25
$ws.Range("L134").Value = 97071.75
$ws.Range("N134").Value = -102141.75
